$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update period / publish-date header labels (rolling one year forward)
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("D9").Value = "1399-01-27 (10)"
$ws.Range("E9").Value = "1400-02-05 (8)"
$ws.Range("F9").Value = "1401-02-10 (9)"
$ws.Range("G9").Value = "1402-02-12 (10)"
$ws.Range("H9").Value = "1402-02-12 (2)"

# Update balance-sheet figures with the refreshed database values
$ws.Range("D12").Value = 11019
$ws.Range("E12").Value = 76213
$ws.Range("F12").Value = 308134
$ws.Range("G12").Value = 490432
$ws.Range("H12").Value = 726000
$ws.Range("G13").Value = 1050000
$ws.Range("H13").Value = 2888500
$ws.Range("D14").Value = 389862
$ws.Range("E14").Value = 515006
$ws.Range("F14").Value = 75801
$ws.Range("G14").Value = 184916
$ws.Range("H14").Value = 246119
$ws.Range("D15").Value = 470121
$ws.Range("E15").Value = 660092
$ws.Range("F15").Value = 821947
$ws.Range("G15").Value = 1798750
$ws.Range("H15").Value = 2960057
$ws.Range("D16").Value = 50856
$ws.Range("E16").Value = 110085
$ws.Range("F16").Value = 341105
$ws.Range("G16").Value = 359945
$ws.Range("H16").Value = 279916
$ws.Range("E17").Value = 1691
$ws.Range("F17").Value = 1580
$ws.Range("H17").Value = 0
$ws.Range("D18").Value = 921858
$ws.Range("E18").Value = 1363087
$ws.Range("F18").Value = 1548567
$ws.Range("G18").Value = 3885623
$ws.Range("H18").Value = 7100592
$ws.Range("D19").Value = 1458
$ws.Range("E19").Value = 7199
$ws.Range("F19").Value = 8248
$ws.Range("G19").Value = 9236
$ws.Range("H19").Value = 12247
$ws.Range("D20").Value = 186714
$ws.Range("E20").Value = 186556
$ws.Range("F20").Value = 805828
$ws.Range("G20").Value = 806204
$ws.Range("H20").Value = 779420
$ws.Range("D22").Value = 1387412
$ws.Range("E22").Value = 1318522
$ws.Range("F22").Value = 1424710
$ws.Range("G22").Value = 1660006
$ws.Range("H22").Value = 2219884
$ws.Range("D23").Value = 13144
$ws.Range("E23").Value = 136757
$ws.Range("F23").Value = 124387
$ws.Range("G23").Value = 99540
$ws.Range("H23").Value = 96051
$ws.Range("D25").Value = 114583
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 51815
$ws.Range("G25").Value = 50490
$ws.Range("H25").Value = 120
$ws.Range("D26").Value = 1703311
$ws.Range("E26").Value = 1649034
$ws.Range("F26").Value = 2414988
$ws.Range("G26").Value = 2625476
$ws.Range("H26").Value = 3107722
$ws.Range("D27").Value = 2625169
$ws.Range("E27").Value = 3012121
$ws.Range("F27").Value = 3963555
$ws.Range("G27").Value = 6511099
$ws.Range("H27").Value = 10208314
$ws.Range("D29").Value = 528581
$ws.Range("E29").Value = 539576
$ws.Range("F29").Value = 515676
$ws.Range("G29").Value = 943258
$ws.Range("H29").Value = 1158148
$ws.Range("D31").Value = 168800
$ws.Range("E31").Value = 160102
$ws.Range("F31").Value = 159231
$ws.Range("G31").Value = 479879
$ws.Range("H31").Value = 862167
$ws.Range("D32").Value = 52232
$ws.Range("E32").Value = 96463
$ws.Range("F32").Value = 199204
$ws.Range("G32").Value = 357012
$ws.Range("H32").Value = 376960
$ws.Range("D33").Value = 531579
$ws.Range("E33").Value = 103292
$ws.Range("F33").Value = 77221
$ws.Range("G33").Value = 106614
$ws.Range("H33").Value = 111167
$ws.Range("D34").Value = 190002
$ws.Range("E34").Value = 158770
$ws.Range("F34").Value = 0
$ws.Range("H34").Value = 1007373
$ws.Range("D37").Value = 1471194
$ws.Range("E37").Value = 1058203
$ws.Range("F37").Value = 951332
$ws.Range("G37").Value = 1886763
$ws.Range("H37").Value = 3515815
$ws.Range("D39").Value = "-"
$ws.Range("D40").Value = 225147
$ws.Range("E40").Value = 123506
$ws.Range("F40").Value = 0
$ws.Range("D41").Value = 190476
$ws.Range("E41").Value = 279753
$ws.Range("F41").Value = 349478
$ws.Range("G41").Value = 1003755
$ws.Range("H41").Value = 1443338
$ws.Range("D42").Value = 415623
$ws.Range("E42").Value = 403259
$ws.Range("F42").Value = 349478
$ws.Range("G42").Value = 1003755
$ws.Range("H42").Value = 1443338
$ws.Range("D43").Value = 1886817
$ws.Range("E43").Value = 1461462
$ws.Range("F43").Value = 1300810
$ws.Range("G43").Value = 2890518
$ws.Range("H43").Value = 4959153
$ws.Range("E45").Value = 1100000
$ws.Range("F48").Value = -57381
$ws.Range("G48").Value = -44954
$ws.Range("H48").Value = -110837
$ws.Range("D49").Value = 0
$ws.Range("F49").Value = 12010
$ws.Range("G49").Value = 4097
$ws.Range("H49").Value = 40000
$ws.Range("E50").Value = 73615
$ws.Range("F50").Value = 110000
$ws.Range("D52").Value = "-"
$ws.Range("D54").Value = "-"
$ws.Range("D56").Value = 133352
$ws.Range("E56").Value = 377044
$ws.Range("F56").Value = 1498116
$ws.Range("G56").Value = 2451438
$ws.Range("H56").Value = 4109998
$ws.Range("D57").Value = 738352
$ws.Range("E57").Value = 1550659
$ws.Range("F57").Value = 2662745
$ws.Range("G57").Value = 3620581
$ws.Range("H57").Value = 5249161
$ws.Range("D58").Value = 2625169
$ws.Range("E58").Value = 3012121
$ws.Range("F58").Value = 3963555
$ws.Range("G58").Value = 6511099
$ws.Range("H58").Value = 10208314
